$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "ScriptType"
$ws.Range("B1").Value = "UI"
$ws.Range("C1").Value = ""
$ws.Range("D1").Value = "Component"
$ws.Range("E1").Value = "URL"
$ws.Range("F1").Value = "alert"
$ws.Range("G1").Value = "confirm"
$ws.Range("H1").Value = "prompt"

# Row 2
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "Xpath"
$ws.Range("E2").Value = "URL"
$ws.Range("F2").Value = "alert"
$ws.Range("G2").Value = "confirm"
$ws.Range("H2").Value = "prompt"

# Row 3
$ws.Range("A3").Value = "TCID"
$ws.Range("B3").Value = "TestCase"
$ws.Range("C3").Value = "TAG"
$ws.Range("D3").Value = "DependsOn"
$ws.Range("E3").Value = "URL"
$ws.Range("F3").Value = "alert"
$ws.Range("G3").Value = "confirm"
$ws.Range("H3").Value = "prompt"
